# Auto-generated edit script applying recalculated profit values per sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 15297.533
$ws.Range("I101").Value = 533
$ws.Range("J101").Value = 22679.8
$ws.Range("K101").Value = 1599
$ws.Range("L101").Value = 68039.39999999999
$ws.Range("M101").Value = 23
$ws.Range("N101").Value = -71283.39999999999
$ws.Range("H132").Value = 3028341.2
$ws.Range("I132").Value = 613383.75
$ws.Range("J132").Value = 27781654
$ws.Range("K132").Value = 1840151.25
$ws.Range("L132").Value = 83344962
$ws.Range("M132").Value = -1837621.25
$ws.Range("N132").Value = -83350022
$ws.Range("H137").Value = 11997742
$ws.Range("I137").Value = 2718128.2
$ws.Range("J137").Value = 50803400
$ws.Range("K137").Value = 8154384.600000001
$ws.Range("L137").Value = 152410200
$ws.Range("M137").Value = -8151834.600000001
$ws.Range("N137").Value = -152415300
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3244136.8
$ws.Range("I32").Value = 4461388.5
$ws.Range("J32").Value = 22000
$ws.Range("K32").Value = 4461388.5
$ws.Range("L32").Value = 22000
$ws.Range("M32").Value = -4461101.5
$ws.Range("N32").Value = -22574
$ws.Range("H37").Value = 9524.941000000001
$ws.Range("J37").Value = 14710
$ws.Range("L37").Value = 14710
$ws.Range("N37").Value = -15256
$ws.Range("H61").Value = 3466742.2
$ws.Range("I61").Value = 1895254
$ws.Range("J61").Value = 8405705
$ws.Range("K61").Value = 1895254
$ws.Range("L61").Value = 8405705
$ws.Range("M61").Value = -1895042
$ws.Range("N61").Value = -8406129
$ws.Range("H132").Value = 35445990
$ws.Range("I132").Value = 44429430
$ws.Range("J132").Value = 9618606
$ws.Range("K132").Value = 133288290
$ws.Range("L132").Value = 28855818
$ws.Range("M132").Value = -133285760
$ws.Range("N132").Value = -28860878
$ws.Range("H136").Value = 3466742.2
$ws.Range("I136").Value = 1895254
$ws.Range("J136").Value = 8405705
$ws.Range("K136").Value = 5685762
$ws.Range("L136").Value = 25217115
$ws.Range("M136").Value = -5683212
$ws.Range("N136").Value = -25222215
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2692474.8
$ws.Range("I31").Value = 5209744
$ws.Range("J31").Value = 7387.6
$ws.Range("K31").Value = 5209744
$ws.Range("L31").Value = 7387.6
$ws.Range("M31").Value = -5209449
$ws.Range("N31").Value = -7977.6
$ws.Range("H34").Value = 2692474.8
$ws.Range("I34").Value = 5209744
$ws.Range("J34").Value = 7387.6
$ws.Range("K34").Value = 5209744
$ws.Range("L34").Value = 7387.6
$ws.Range("M34").Value = -5209542
$ws.Range("N34").Value = -7791.6
$ws.Range("H51").Value = 30543.889
$ws.Range("J51").Value = 10699.286
$ws.Range("L51").Value = 10699.286
$ws.Range("N51").Value = -12171.286
$ws.Range("H58").Value = 1981971
$ws.Range("I58").Value = 7677.467
$ws.Range("J58").Value = 5683771.5
$ws.Range("K58").Value = 7677.467
$ws.Range("L58").Value = 5683771.5
$ws.Range("M58").Value = -7474.467
$ws.Range("N58").Value = -5684177.5
$ws.Range("H60").Value = 21801.75
$ws.Range("I60").Value = 53900
$ws.Range("J60").Value = 11102.333
$ws.Range("K60").Value = 53900
$ws.Range("L60").Value = 11102.333
$ws.Range("M60").Value = -53389
$ws.Range("N60").Value = -12124.333
$ws.Range("H61").Value = 30543.889
$ws.Range("J61").Value = 10699.286
$ws.Range("L61").Value = 10699.286
$ws.Range("N61").Value = -11395.286
$ws.Range("H68").Value = 18120
$ws.Range("J68").Value = 18120
$ws.Range("L68").Value = 18120
$ws.Range("N68").Value = -19618
$ws.Range("H71").Value = 18120
$ws.Range("J71").Value = 18120
$ws.Range("L71").Value = 54360
$ws.Range("N71").Value = -61848
$ws.Range("H74").Value = 17380.637
$ws.Range("J74").Value = 18990.2
$ws.Range("L74").Value = 18990.2
$ws.Range("N74").Value = -20738.2
$ws.Range("H77").Value = 17380.637
$ws.Range("J77").Value = 18990.2
$ws.Range("L77").Value = 56970.60000000001
$ws.Range("N77").Value = -65706.60000000001
$ws.Range("H134").Value = 1002268.1
$ws.Range("I134").Value = 1578.5927
$ws.Range("J134").Value = 3080623.2
$ws.Range("K134").Value = 4735.7781
$ws.Range("L134").Value = 9241869.600000001
$ws.Range("M134").Value = -2200.7781
$ws.Range("N134").Value = -9246939.600000001
$ws.Range("H136").Value = 1981971
$ws.Range("I136").Value = 7677.467
$ws.Range("J136").Value = 5683771.5
$ws.Range("K136").Value = 23032.401
$ws.Range("L136").Value = 17051314.5
$ws.Range("M136").Value = -20482.401
$ws.Range("N136").Value = -17056414.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8887.066000000001
$ws.Range("I3").Value = 7358.8335
$ws.Range("K3").Value = 22076.5005
$ws.Range("M3").Value = -21964.5005
$ws.Range("H108").Value = 1885.5834
$ws.Range("I108").Value = 325.4
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 976.1999999999999
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = 1903.8
$ws.Range("N108").Value = -14760
$ws.Range("H118").Value = 2234.3
$ws.Range("I118").Value = 1468.4286
$ws.Range("K118").Value = 4405.2858
$ws.Range("M118").Value = -3162.2858
$ws.Range("H132").Value = 1676.625
$ws.Range("I132").Value = 1298.3334
$ws.Range("J132").Value = 1743.3823
$ws.Range("K132").Value = 11685.0006
$ws.Range("L132").Value = 15690.4407
$ws.Range("M132").Value = -9155.000599999999
$ws.Range("N132").Value = -20750.4407
$ws.Range("H133").Value = 3095.7576
$ws.Range("I133").Value = 3243.077
$ws.Range("J133").Value = 3000
$ws.Range("K133").Value = 9729.231
$ws.Range("L133").Value = 9000
$ws.Range("M133").Value = -4669.231
$ws.Range("N133").Value = -19120
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 15390162
$ws.Range("I132").Value = 22512968
$ws.Range("J132").Value = 8267355
$ws.Range("K132").Value = 67538904
$ws.Range("L132").Value = 24802065
$ws.Range("M132").Value = -67536374
$ws.Range("N132").Value = -24807125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 22729010
$ws.Range("I22").Value = 507.7143
$ws.Range("J22").Value = 33335644
$ws.Range("K22").Value = 507.7143
$ws.Range("L22").Value = 33335644
$ws.Range("M22").Value = -212.7143
$ws.Range("N22").Value = -33336234
$ws.Range("H27").Value = 22729010
$ws.Range("I27").Value = 507.7143
$ws.Range("J27").Value = 33335644
$ws.Range("K27").Value = 507.7143
$ws.Range("L27").Value = 33335644
$ws.Range("M27").Value = -400.7143
$ws.Range("N27").Value = -33335858
$ws.Range("H55").Value = 35714496
$ws.Range("I55").Value = 50000156
$ws.Range("J55").Value = 335
$ws.Range("K55").Value = 50000156
$ws.Range("L55").Value = 335
$ws.Range("M55").Value = -49999983
$ws.Range("N55").Value = -681
$ws.Range("H122").Value = 12295193
$ws.Range("I122").Value = 1420697.6
$ws.Range("J122").Value = 66667670
$ws.Range("K122").Value = 4262092.800000001
$ws.Range("L122").Value = 200003010
$ws.Range("M122").Value = -4259642.800000001
$ws.Range("N122").Value = -200007910
$ws.Range("H132").Value = 6502601
$ws.Range("I132").Value = 11914843
$ws.Range("J132").Value = 7910.4
$ws.Range("K132").Value = 35744529
$ws.Range("L132").Value = 23731.2
$ws.Range("M132").Value = -35741999
$ws.Range("N132").Value = -28791.2
$ws.Range("H135").Value = 39056.5
$ws.Range("J135").Value = 39056.5
$ws.Range("L135").Value = 39056.5
$ws.Range("N135").Value = -49196.5
$ws.Range("H136").Value = 2453171.8
$ws.Range("I136").Value = 2606230.5
$ws.Range("K136").Value = 7818691.5
$ws.Range("M136").Value = -7816141.5
$ws.Range("H139").Value = 48605.26
$ws.Range("J139").Value = 48605.26
$ws.Range("L139").Value = 48605.26
$ws.Range("N139").Value = -58885.26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1591940.8
$ws.Range("I132").Value = 5131.7617
$ws.Range("J132").Value = 5294495
$ws.Range("K132").Value = 15395.2851
$ws.Range("L132").Value = 15883485
$ws.Range("M132").Value = -12865.2851
$ws.Range("N132").Value = -15888545
$ws.Range("H136").Value = 1647.5333
$ws.Range("I136").Value = 1373
$ws.Range("J136").Value = 4118.3335
$ws.Range("K136").Value = 4119
$ws.Range("L136").Value = 12355.0005
$ws.Range("M136").Value = -1569
$ws.Range("N136").Value = -17455.0005

Write-Host "Applied all updates"
